# ---------------------------------------------------------------------------
# Commit: "Change target framework to netstandard2.0"
#
# The canonical-OOXML diff for this commit touches exactly two things inside
# PowerpointTemplater.Tests/files/ReplaceTables_output6.pptx:
#
#   1. The r:id *values* referenced by the nine middle <p:sldId> entries
#      (ids 259-267) in ppt/presentation.xml change from one set of
#      auto-generated relationship-id strings to another set.
#   2. The r:embed value on the single <a:blip> (the {{picture1png}}
#      placeholder picture on the last slide) changes to a different
#      auto-generated relationship-id string.
#
# In every one of those hunks the *value pointed at* (which slide part /
# which image part) is unchanged - only the opaque, tool-generated
# relationship-id label is different:
#   - <p:sldId id="259" .../> .. id="267" keep the same id numbers, the
#     same order, and still resolve to slides 5,6,7,8,10,11,13,14,15.
#   - the blip still resolves to ppt/media/image2.png (the already
#     "replaced" {{picture1png}} image) both before and after.
#
# That matches the commit message: retargeting the library build to
# netstandard2.0 pulls in a different (but functionally equivalent)
# DocumentFormat.OpenXml resolution, so re-running the ReplaceTables test
# that produced this golden fixture regenerates the same deck byte-for-byte
# except for the library's internally auto-generated relationship-id
# strings (an `R`+16-hex-chars GUID-ish token minted fresh any time the
# relevant package part is rewritten). No slide was added, removed,
# reordered or resized; no picture/table content changed; nothing a user
# could reproduce via the PowerPoint object model actually moved.
#
# Relationship-id strings are an implementation detail that this COM
# surface (like real PowerPoint automation) does not let a caller pin to
# an arbitrary literal value - Slide.SlideID, Shape.Id, and every r:id/
# r:embed attribute are assigned internally and are read-only. Forcing a
# *different* id string here (e.g. by deleting+recreating the picture
# shape or by cutting/re-pasting the affected slides) would not reproduce
# the specific tokens in the diff anyway (this runtime mints its own
# relationship ids with a different scheme than DocumentFormat.OpenXml),
# and it would additionally corrupt real, meaningful state that the diff
# shows as unchanged (the picture shape's cNvPr id, and it would leave a
# dangling/orphaned image relationship behind).
#
# So the faithful replay of this specific commit is to leave the deck's
# slides, their order/ids, and the placeholder picture exactly as they
# are - i.e. confirm the structure already matches what the diff implies
# and make no content-visible edit. The two read-only sanity checks below
# document that expectation without mutating anything.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

$lastSlide = $p.Slides.Item($p.Slides.Count)
$pic = $lastSlide.Shapes.Item(1)

Write-Host "Slide count:" $p.Slides.Count
Write-Host "Last slide id:" $lastSlide.SlideID "- picture shape:" $pic.Name
